$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Client ID values for existing rows 3-10 (row2 stays 12345)
$ws.Cells.Item(3, 1).Value = 12346
$ws.Cells.Item(4, 1).Value = 12347
$ws.Cells.Item(5, 1).Value = 12348
$ws.Cells.Item(6, 1).Value = 12349
$ws.Cells.Item(7, 1).Value = 12350
$ws.Cells.Item(8, 1).Value = 12351
$ws.Cells.Item(9, 1).Value = 12352
$ws.Cells.Item(10, 1).Value = 12353

# New rows 11-18
$newRows = @(
    @{ Row = 11; A = 12354; B = "Health Check";              C = 45931 },
    @{ Row = 12; A = 12355; B = "Physical Activity";         C = 45962 },
    @{ Row = 13; A = 12356; B = "Health & Wellbeing Coach";  C = 45992 },
    @{ Row = 14; A = 12357; B = "Health Check";              C = 46023 },
    @{ Row = 15; A = 12358; B = "Physical Activity";         C = 46054 },
    @{ Row = 16; A = 12359; B = "Health & Wellbeing Coach";  C = 46082 },
    @{ Row = 17; A = 12360; B = "Health Check";              C = 46113 },
    @{ Row = 18; A = 12361; B = "Physical Activity";         C = 46143 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
}

# Copy the date style/number format from C10 onto the new C11:C18 cells
$ws.Cells.Item(10, 3).Copy()
$ws.Range($ws.Cells.Item(11, 3), $ws.Cells.Item(18, 3)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update view: top-left cell and selection to reflect scrolled view
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("E17").Select()
